# Add season record columns (Wins, Losses, Ties) to the player stats sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): labels for the new columns, matching the
#     existing header style (bold font, thin border, centered) ---
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-54): every player on the roster shares the team's
#     season record: 85 wins, 77 losses, 0 ties ---
$lastRow = 54
$ws.Range("AD2:AD" + $lastRow).Value = 85
$ws.Range("AE2:AE" + $lastRow).Value = 77
$ws.Range("AF2:AF" + $lastRow).Value = 0

Write-Output "Added Wins/Losses/Ties columns"
